# Append a new airdrop-list row (row 4) to the single worksheet, mirroring
# the format/styles already used by the existing rows, and select the new
# full used range (matches the author re-saving the file after adding a row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row values -------------------------------------------------------
$ws.Range("A4").Value = 43133.096827175927
$ws.Range("B4").Value = "@Larsblm"
$ws.Range("C4").Value = "https://twitter.com/Larsblm/status/959369283116716036"
$ws.Range("D4").Value = "0x645a06e738adea003014583f5e42508e6ecbb040"

# --- copy row 3's formatting down to row 4 (style/number format/height) --
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(3).RowHeight

# --- hyperlink for the retweet-link cell, like C2/C3 ----------------------
$ws.Hyperlinks.Add($ws.Range("C4"), "https://twitter.com/Larsblm/status/959369283116716036")

# re-apply the plain style used by the other retweet-link cells, undoing
# the built-in "Hyperlink" style that Hyperlinks.Add stamps onto the cell
$ws.Range("C4").Style = $ws.Range("C3").Style
$ws.Range("C3:C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C4").Value = "https://twitter.com/Larsblm/status/959369283116716036"

# --- keep the full-sheet selection in sync with the new used range -------
$ws.Range("A1:XFD4").Select()
